$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing 2021 column (R) into the new 2022
# column (S) for the header row and every data row, then fill in the
# new 2022 values.
$ws.Range("R4:R14").Copy()
$ws.Range("S4:S14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("S4").Value = 2022
$ws.Range("S5").Value = 99.5
$ws.Range("S6").Value = 99.358544044156048
$ws.Range("S7").Value = 99.400057479522914
$ws.Range("S8").Value = 99.513194978221875
$ws.Range("S9").Value = 99.232429839290006
$ws.Range("S10").Value = 99.453093666824671
$ws.Range("S11").Value = 99.686258104998956
$ws.Range("S12").Value = 99.42525365081228
$ws.Range("S13").Value = 99.561275226674468
$ws.Range("S14").Value = 99.831561216970215

# Match the author's final selection.
$ws.Range("U6").Select()
